$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$tr = $s.Shapes.Item(1).TextFrame.TextRange
$tr.Characters(1, $tr.Length).Text = "Below section-level"
